$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4999
$ws.Range("I74").Value = 4999
$ws.Range("K74").Value = 4999
$ws.Range("M74").Value = -4063
$ws.Range("H77").Value = 4999
$ws.Range("I77").Value = 4999
$ws.Range("K77").Value = 24995
$ws.Range("M77").Value = -20315
$ws.Range("H96").Value = 3533.3333
$ws.Range("I96").Value = 5000
$ws.Range("J96").Value = 600
$ws.Range("K96").Value = 15000
$ws.Range("L96").Value = 1800
$ws.Range("M96").Value = -13627
$ws.Range("N96").Value = -4546
$ws.Range("H100").Value = 4833.3335
$ws.Range("I100").Value = 4833.3335
$ws.Range("K100").Value = 4833.3335
$ws.Range("M100").Value = -4292.3335
$ws.Range("H125").Value = 2666.3333
$ws.Range("I125").Value = 2500
$ws.Range("J125").Value = 2999
$ws.Range("K125").Value = 22500
$ws.Range("L125").Value = 26991
$ws.Range("M125").Value = -20040
$ws.Range("N125").Value = -31911
$ws.Range("H132").Value = 2743.9092
$ws.Range("I132").Value = 1740.4286
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 5221.2858
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -2691.2858
$ws.Range("N132").Value = -18560
$ws.Range("H138").Value = 3788.5264
$ws.Range("I138").Value = 1999.5
$ws.Range("K138").Value = 5998.5
$ws.Range("M138").Value = -858.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4413.5386
$ws.Range("I32").Value = 3107.818
$ws.Range("J32").Value = 11595
$ws.Range("K32").Value = 3107.818
$ws.Range("L32").Value = 11595
$ws.Range("M32").Value = -2820.818
$ws.Range("N32").Value = -12169
$ws.Range("H74").Value = 3528.25
$ws.Range("I74").Value = 3528.25
$ws.Range("K74").Value = 3528.25
$ws.Range("M74").Value = -2654.25
$ws.Range("H77").Value = 3528.25
$ws.Range("I77").Value = 3528.25
$ws.Range("K77").Value = 17641.25
$ws.Range("M77").Value = -13273.25
$ws.Range("H92").Value = 18583.334
$ws.Range("J92").Value = 25375
$ws.Range("L92").Value = 25375
$ws.Range("N92").Value = -30367
$ws.Range("H119").Value = 99900
$ws.Range("J119").Value = 99900
$ws.Range("L119").Value = 99900
$ws.Range("N119").Value = -109576
$ws.Range("H122").Value = 16650
$ws.Range("I122").Value = 14975
$ws.Range("J122").Value = 20000
$ws.Range("K122").Value = 44925
$ws.Range("L122").Value = 60000
$ws.Range("M122").Value = -42475
$ws.Range("N122").Value = -64900
$ws.Range("H132").Value = 3100.8333
$ws.Range("I132").Value = 3100.8333
$ws.Range("K132").Value = 9302.499899999999
$ws.Range("M132").Value = -6772.499899999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 776.6667
$ws.Range("I80").Value = 679.4
$ws.Range("J80").Value = 898.25
$ws.Range("K80").Value = 679.4
$ws.Range("L80").Value = 898.25
$ws.Range("M80").Value = 318.6
$ws.Range("N80").Value = -2894.25
$ws.Range("H83").Value = 776.6667
$ws.Range("I83").Value = 679.4
$ws.Range("J83").Value = 898.25
$ws.Range("K83").Value = 3397
$ws.Range("L83").Value = 4491.25
$ws.Range("M83").Value = 1595
$ws.Range("N83").Value = -14475.25
$ws.Range("H96").Value = 8900
$ws.Range("I96").Value = 8900
$ws.Range("K96").Value = 8900
$ws.Range("M96").Value = -6154
$ws.Range("H107").Value = 13393
$ws.Range("I107").Value = 6618.857
$ws.Range("J107").Value = 29199.334
$ws.Range("K107").Value = 6618.857
$ws.Range("L107").Value = 29199.334
$ws.Range("M107").Value = -4698.857
$ws.Range("N107").Value = -33039.334
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 9999.5
$ws.Range("I86").Value = 3999
$ws.Range("K86").Value = 3999
$ws.Range("M86").Value = -2876
$ws.Range("H89").Value = 9999.5
$ws.Range("I89").Value = 3999
$ws.Range("K89").Value = 19995
$ws.Range("M89").Value = -14379
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1700.875
$ws.Range("J121").Value = 2376.4
$ws.Range("L121").Value = 7129.200000000001
$ws.Range("N121").Value = -9749.200000000001
$ws.Range("H133").Value = 4000
$ws.Range("I133").Value = 4000
$ws.Range("K133").Value = 12000
$ws.Range("M133").Value = -6940
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4360.9
$ws.Range("I132").Value = 4657.8887
$ws.Range("J132").Value = 1688
$ws.Range("K132").Value = 13973.6661
$ws.Range("L132").Value = 5064
$ws.Range("M132").Value = -11443.6661
$ws.Range("N132").Value = -10124
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 900
$ws.Range("I46").Value = 800
$ws.Range("K46").Value = 800
$ws.Range("M46").Value = -612
$ws.Range("H74").Value = 100000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 100000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 100000
$ws.Range("M74").Value = ""
$ws.Range("N74").Value = -101996
$ws.Range("H77").Value = 100000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 100000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 300000
$ws.Range("M77").Value = ""
$ws.Range("N77").Value = -309984
$ws.Range("H122").Value = 7448
$ws.Range("I122").Value = 7448
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 22344
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -19894
$ws.Range("N122").Value = ""
$ws.Range("H132").Value = 16276.615
$ws.Range("I132").Value = 16966.334
$ws.Range("K132").Value = 50899.00199999999
$ws.Range("M132").Value = -48369.00199999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8545.416999999999
$ws.Range("J62").Value = 3511
$ws.Range("L62").Value = 3511
$ws.Range("N62").Value = -4759
$ws.Range("H65").Value = 8545.416999999999
$ws.Range("J65").Value = 3511
$ws.Range("L65").Value = 17555
$ws.Range("N65").Value = -23795
$ws.Range("H107").Value = 525.1539
$ws.Range("I107").Value = 402.8
$ws.Range("K107").Value = 1208.4
$ws.Range("M107").Value = 711.5999999999999
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""
$ws.Range("H122").Value = 1921
$ws.Range("I122").Value = 1683.4286
$ws.Range("J122").Value = 2752.5
$ws.Range("K122").Value = 5050.2858
$ws.Range("L122").Value = 8257.5
$ws.Range("M122").Value = -2600.2858
$ws.Range("N122").Value = -13157.5
$ws.Range("H132").Value = 3667.3333
$ws.Range("I132").Value = 3667.3333
$ws.Range("K132").Value = 11001.9999
$ws.Range("M132").Value = -8471.999899999999
